# Update visitor/attendance counts (column F) on the "展览" and "全部类型"
# sheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 586
$ws1.Range("F5").Value  = 1096
$ws1.Range("F6").Value  = 14085
$ws1.Range("F7").Value  = 15591
$ws1.Range("F8").Value  = 8
$ws1.Range("F9").Value  = 51
$ws1.Range("F20").Value = 1215
$ws1.Range("F23").Value = 6091
$ws1.Range("F25").Value = 1086
$ws1.Range("F26").Value = 5545
$ws1.Range("F27").Value = 72
$ws1.Range("F30").Value = 4510

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 586
$ws4.Range("F5").Value  = 1096
$ws4.Range("F6").Value  = 14085
$ws4.Range("F7").Value  = 15591
$ws4.Range("F8").Value  = 8
$ws4.Range("F9").Value  = 51
$ws4.Range("F20").Value = 1215
$ws4.Range("F24").Value = 6091
$ws4.Range("F26").Value = 1086
$ws4.Range("F27").Value = 5545
$ws4.Range("F28").Value = 72
$ws4.Range("F31").Value = 4510
